$d = $word.ActiveDocument

# --- Step 1: Bold the first three list items ---
# "Lugares reservados", "Lugares Vazios", "Lugares que o comboio tem"
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Bold = 1

$p2 = $d.Paragraphs.Item(2)
$p2.Range.Bold = 1

$p3 = $d.Paragraphs.Item(3)
$p3.Range.Bold = 1

# --- Step 2: Move the _GoBack bookmark from the last paragraph to the end of
#     paragraph 3's text (after "Lugares que o comboio tem", still inside w:p 3) ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $gb = $d.Bookmarks("_GoBack")
    $gb.Delete()
}

# Refetch paragraph 3 (formatting calls above should not have changed paragraph count)
$p3 = $d.Paragraphs.Item(3)
$insertPos = $p3.Range.End - 1

# Work around a quirk where adding a bookmark collapsed exactly on the gap
# immediately before a paragraph mark can misplace it: temporarily insert a
# marker character after the text, add the bookmark there, then remove the
# marker again. The bookmark stays anchored at the correct position.
$marker = $d.Range($insertPos, $insertPos)
$marker.InsertAfter("@")
$bmRange = $d.Range($insertPos, $insertPos)
$bmRange.Bookmarks.Add("_GoBack")
$d.Range($insertPos, $insertPos + 1).Delete()

# --- Step 3: Insert a blank paragraph right after paragraph 3 ---
$p3 = $d.Paragraphs.Item(3)
$afterP3 = $d.Range($p3.Range.End, $p3.Range.End)
[void]$afterP3.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>")

# --- Step 4: Merge the multi-run paragraphs into a single run each ---
function Find-ParagraphByPrefix($doc, $prefix) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text
        if ($t.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

$target1 = Find-ParagraphByPrefix $d "Viagens internacionais"
$rng1 = $d.Range($target1.Range.Start, $target1.Range.End - 1)
$rng1.Text = "__TEMP_PLACEHOLDER_1__"
$target1b = Find-ParagraphByPrefix $d "__TEMP_PLACEHOLDER_1__"
$rng1b = $d.Range($target1b.Range.Start, $target1b.Range.End - 1)
$rng1b.Text = "Viagens internacionais feitas antes/depois do dia X"

$target2 = Find-ParagraphByPrefix $d "Viagem internacionais"
$rng2 = $d.Range($target2.Range.Start, $target2.Range.End - 1)
$rng2.Text = "__TEMP_PLACEHOLDER_2__"
$target2b = Find-ParagraphByPrefix $d "__TEMP_PLACEHOLDER_2__"
$rng2b = $d.Range($target2b.Range.Start, $target2b.Range.End - 1)
$rng2b.Text = "Viagem internacionais feitas num dia/mês/ano"

Write-Host "done"
